$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: temporarily format column D cells as Text so that numeric-looking
# strings (e.g. "6.78", "1.00") are not auto-converted to numbers on write.
$ws.Range("D2:D51").NumberFormat = "@"

# Step 2: write the updated Price (D) and Volume(1h) (E) values.
$ws.Range("D2").Value = "58.777.74"
$ws.Range("E2").Value = "  -1.11%  "
$ws.Range("D3").Value = "2.595.21"
$ws.Range("E3").Value = "  -1.56%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "553.65"
$ws.Range("E5").Value = "  +3.37%  "
$ws.Range("D6").Value = "143.66"
$ws.Range("E6").Value = "  -0.58%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "0.598"
$ws.Range("E8").Value = "  +4.87%  "
$ws.Range("D9").Value = "6.78"
$ws.Range("E9").Value = "  +2.31%  "
$ws.Range("E10").Value = "  -1.21%  "
$ws.Range("E11").Value = "  +4.85%  "
$ws.Range("D12").Value = "0.335"
$ws.Range("E12").Value = "  -0.42%  "
$ws.Range("D13").Value = "3.057.36"
$ws.Range("E13").Value = "  -1.68%  "
$ws.Range("D14").Value = "59.075.39"
$ws.Range("E14").Value = "  -0.48%  "
$ws.Range("D15").Value = "20.85"
$ws.Range("E15").Value = "  -0.98%  "
$ws.Range("D16").Value = "2.616.61"
$ws.Range("E16").Value = "  -1.70%  "
$ws.Range("E17").Value = "  -1.96%  "
$ws.Range("D18").Value = "4.44"
$ws.Range("E18").Value = "  +1.40%  "
$ws.Range("D19").Value = "337.41"
$ws.Range("E19").Value = "  -0.54%  "
$ws.Range("D20").Value = "10.05"
$ws.Range("E20").Value = "  -2.65%  "
$ws.Range("D21").Value = "6.14"
$ws.Range("E21").Value = "  -2.38%  "
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("D23").Value = "66.42"
$ws.Range("E23").Value = "  -0.93%  "
$ws.Range("D24").Value = "0.425"
$ws.Range("E24").Value = "  +2.57%  "
$ws.Range("E25").Value = "  -0.40%  "
$ws.Range("E26").Value = "  -3.52%  "
$ws.Range("D27").Value = "7.12"
$ws.Range("E27").Value = "  -1.84%  "
$ws.Range("D28").Value = "0.0₃0758"
$ws.Range("E28").Value = "  +2.19%  "
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("E30").Value = "  +1.43%  "
$ws.Range("D31").Value = "5.93"
$ws.Range("E31").Value = "  +1.88%  "
$ws.Range("D32").Value = "154.45"
$ws.Range("E32").Value = "  +2.01%  "
$ws.Range("E33").Value = "  +0.40%  "
$ws.Range("D34").Value = "3.91"
$ws.Range("E34").Value = "  -1.80%  "
$ws.Range("D35").Value = "0.874"
$ws.Range("E35").Value = "  +4.04%  "
$ws.Range("D36").Value = "37.15"
$ws.Range("E36").Value = "  +0.18%  "
$ws.Range("E37").Value = "  -1.57%  "
$ws.Range("D38").Value = "1.45"
$ws.Range("E38").Value = "  +0.64%  "
$ws.Range("D39").Value = "0.822"
$ws.Range("E39").Value = "  -1.08%  "
$ws.Range("D40").Value = "3.60"
$ws.Range("E40").Value = "  +0.74%  "
$ws.Range("D41").Value = "281.04"
$ws.Range("E41").Value = "  -2.45%  "
$ws.Range("D42").Value = "0.997"
$ws.Range("E42").Value = "  -0.26%  "
$ws.Range("D43").Value = "0.596"
$ws.Range("E43").Value = "  -1.16%  "
$ws.Range("D44").Value = "0.0952"
$ws.Range("E44").Value = "  +0.69%  "
$ws.Range("E45").Value = "  -0.78%  "
$ws.Range("D46").Value = "0.0531"
$ws.Range("E46").Value = "  -0.26%  "
$ws.Range("E47").Value = "  +0.78%  "
$ws.Range("D48").Value = "1.919.67"
$ws.Range("E48").Value = "  -2.48%  "
$ws.Range("D49").Value = "4.44"
$ws.Range("E49").Value = "  -2.45%  "
$ws.Range("D50").Value = "17.82"
$ws.Range("E50").Value = "  -2.24%  "
$ws.Range("D51").Value = "114.61"
$ws.Range("E51").Value = "  +3.15%  "

# Step 3: restore the original (default) cell style on column D so no
# stray number-format style is left behind on cells that did not have one.
$ws.Range("D2:D51").Style = $ws.Range("B2").Style

